$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text for the new "NOTES" help row (appended to sharedStrings, becomes row 26)
$notes = @'
bn:NOTES:
Anything between single quotes ' will show in the file name only when ON.
Example: 'REC ~batch'
Anything between double quotes " will show in the file name only when OFF. 
Example: "~operator"
For backward compatibility, when the Prefix field is text only the date and time are appended to the file name.
Example: 'Autosave' will result in file name 'Autosave_20-01-13_1705'.
To show only the text place a single '!' at the start of the Prefix field
Example: '!Autosave' will result in file name 'Autosave'.
To maintain cross platform compatibility, file names may contain only letters, numbers, spaces, 
and the following special characters:  
_ - . ( )
'@

# Write the new cell value in row 26, column A
$ws.Range("A26").Value = $notes

# New style: wrap text (adds a new cellXfs entry) and grow the row to the
# maximum row height so the whole note is visible.
$ws.Range("A26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 409.5

# Make sure the worksheet shows the plain/default view (no frozen/scrolled
# top-left cell, no lingering selection on A2).
[void]$ws.Range("A1").Select()

Write-Host "Appended autosave NOTES row to sheet1"
